$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 06:15:04'
$ws1.Range("A3").Value = 'Total filas: 61'

$data1 = @{
    28 = @('06:15:04','06:15','225_HARAS DEL SUR',0,'LP1912')
    29 = @('05:18:56','06:20','26_HERNANDEZ',62,'LP1912')
    30 = @('04:40:33','06:21','26_HERNANDEZ',101,'LP1912')
    31 = @('05:18:56','06:26','23_HERNANDEZ',68,'LP1912')
    32 = @('04:40:33','06:27','23_HERNANDEZ',107,'LP1912')
    33 = @('06:15:04','06:28','23_HERNANDEZ',13,'LP1912')
    34 = @('04:40:33','06:29','86_EST CHICA-ESC AGRARIA',109,'LP1912')
    35 = @('04:54:03','06:30','86_EST CHICA-ESC AGRARIA',96,'LP1912')
    36 = @('04:40:33','06:31','16_SANTA ANA',111,'LP1912')
    37 = @('05:18:56','06:43','225_C ROCA-H SUR',85,'LP1912')
    38 = @('04:54:03','06:44','225_C ROCA-H SUR',110,'LP1912')
    39 = @('05:18:56','06:46','215C_EL PATO',88,'LP1912')
    40 = @('04:54:03','06:47','215C_EL PATO',113,'LP1912')
    41 = @('05:18:56','06:58','10_OLMOS',100,'LP1912')
    42 = @('05:18:56','06:59','14_ABASTO',101,'LP1912')
    43 = @('06:15:04','07:00','14_ABASTO',45,'LP1912')
    44 = @('06:15:04','07:01','16_SANTA ANA',46,'LP1912')
    45 = @('05:18:56','07:04','15_ABASTO',106,'LP1912')
    46 = @('05:49:40','07:04','23_HERNANDEZ',75,'LP1912')
    47 = @('05:49:40','07:05','15_ABASTO',76,'LP1912')
    48 = @('05:18:56','07:06','225_GOMEZ',108,'LP1912')
    49 = @('05:49:40','07:07','225_GOMEZ',78,'LP1912')
    50 = @('05:18:56','07:11','215A_EL PATO',113,'LP1912')
    51 = @('06:15:04','07:12','215A_EL PATO',57,'LP1912')
    52 = @('05:18:56','07:15','11_ETCHEVERRY',117,'LP1912')
    53 = @('06:15:04','07:16','11_ETCHEVERRY',61,'LP1912')
    54 = @('05:49:40','07:21','26_HERNANDEZ',92,'LP1912')
    55 = @('06:15:04','07:23','10_OLMOS',68,'LP1912')
    56 = @('05:49:40','07:29','10_OLMOS',100,'LP1912')
    57 = @('05:49:40','07:31','11_ETCHEVERRY',102,'LP1912')
    58 = @('05:49:40','07:32','84_COLONIA URQUIZA-ESC 49',103,'LP1912')
    59 = @('06:15:04','07:32','11_ETCHEVERRY',77,'LP1912')
    60 = @('05:49:40','07:36','27_EL RETIRO',107,'LP1912')
    61 = @('06:15:04','07:37','27_EL RETIRO',82,'LP1912')
    62 = @('05:49:40','07:39','10_OLMOS',110,'LP1912')
    63 = @('06:15:04','07:48','14_ABASTO',93,'LP1912')
    64 = @('06:15:04','07:52','215D_EL PATO',97,'LP1912')
    65 = @('06:15:04','08:01','23_HERNANDEZ',106,'LP1912')
    66 = @('06:15:04','08:12','15_ABASTO',117,'LP1912')
}
foreach ($r in $data1.Keys) {
    $row = $data1[$r]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 06:15:04'
$ws2.Range("A3").Value = 'Total filas: 11'

$data2 = @{
    15 = @('06:15:04','07:12','215A_EL PATO',57,'LP1912')
    16 = @('06:15:04','07:52','215D_EL PATO',97,'LP1912')
}
foreach ($r in $data2.Keys) {
    $row = $data2[$r]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
}

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 06:15:04'
$ws3.Range("A3").Value = 'Total filas: 11'

$data3 = @{
    13 = @('06:15:04','07:00','215B_LP-P MOR-1 Y 57',45,'L6173')
    14 = @('05:49:40','07:07','215B_LP-P MOR-1 Y 57',78,'L6173')
    15 = @('05:49:40','07:35','215A_LA PLATA',106,'L6173')
    16 = @('06:15:04','08:07','215C_LA PLATA',112,'L6203')
}
foreach ($r in $data3.Keys) {
    $row = $data3[$r]
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
}
